# Apply the changes described by the commit:
#  - Rename the "Logic" sheet to "Data" (formulas / defined names referring
#    to Logic! auto-update to Data! through the rename).
#  - Update the Results sheet sensitivity test-case row 11 (E11 15 -> 16)
#    and its cached simulation output (Q11), which ripples into the
#    Sensitivity sheet's dependent formulas (B6/G6).
#  - Update the Results sheet selection to match the new state (A2:R12).

$wb = $excel.ActiveWorkbook

# --- Rename "Logic" -> "Data" ------------------------------------------
$logicSheet = $wb.Worksheets.Item("Logic")
$logicSheet.Name = "Data"

# --- Update Results sheet data / simulation output ----------------------
$results = $wb.Worksheets.Item("Results")
$results.Range("E11").Value = 16
$results.Range("Q11").Value = 0.024724261933992901

# Update selection on the Results sheet to A2:R12 (active cell A2)
$null = $results.Range("A2:R12").Select()

# Restore the originally-active sheet/tab ("Data", formerly "Logic") so the
# workbook's active-tab state is unchanged by the Results selection above.
$dataSheet = $wb.Worksheets.Item("Data")
$dataSheet.Activate()
$null = $dataSheet.Range("I14").Select()
